$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("run_1")
$ws.Range("F2").Value = 30.6236538887024
$ws.Range("F3").Value = 30.19364643096924
$ws.Range("F4").Value = 30.04372143745422
$ws.Range("F5").Value = 30.18267607688904
$ws.Range("F6").Value = 30.10561490058899
$ws.Range("F7").Value = 30.17365622520447
$ws.Range("F8").Value = 30.17056465148925
$ws.Range("F9").Value = 30.20643591880798
$ws.Range("F10").Value = 29.94030141830444
$ws.Range("F11").Value = 30.48698830604553
$ws.Range("F12").Value = 29.7547173500061
$ws.Range("F13").Value = 29.50939273834229
$ws.Range("F14").Value = 29.48709893226624
$ws.Range("F15").Value = 29.55073046684265
$ws.Range("F16").Value = 30.92357540130615
$ws.Range("F17").Value = 31.18457174301147
$ws.Range("F18").Value = 30.46803021430969
$ws.Range("F19").Value = 30.0278742313385
$ws.Range("F20").Value = 29.9869487285614
$ws.Range("F21").Value = 30.31540775299072

$ws = $wb.Worksheets.Item("run_2")
$ws.Range("F2").Value = 30.26604509353638
$ws.Range("F3").Value = 30.0528838634491
$ws.Range("F4").Value = 30.05424571037292
$ws.Range("F5").Value = 30.03125309944153
$ws.Range("F6").Value = 29.99723529815674
$ws.Range("F7").Value = 30.04426193237305
$ws.Range("F8").Value = 30.06783390045166
$ws.Range("F9").Value = 30.04636740684509
$ws.Range("F10").Value = 29.96776366233826
$ws.Range("F11").Value = 30.32413291931152
$ws.Range("F12").Value = 29.53144264221192
$ws.Range("F13").Value = 29.52733945846558
$ws.Range("F14").Value = 29.41553902626038
$ws.Range("F15").Value = 29.52752542495728
$ws.Range("F16").Value = 30.85455441474915
$ws.Range("F17").Value = 30.97779631614685
$ws.Range("F18").Value = 30.24300575256348
$ws.Range("F19").Value = 30.25027918815613
$ws.Range("F20").Value = 29.99231290817261
$ws.Range("F21").Value = 30.2535879611969

$ws = $wb.Worksheets.Item("run_3")
$ws.Range("F2").Value = 30.23923921585083
$ws.Range("F3").Value = 30.0996105670929
$ws.Range("F4").Value = 29.98058676719665
$ws.Range("F5").Value = 30.19514441490173
$ws.Range("F6").Value = 30.02770137786865
$ws.Range("F7").Value = 30.17364478111267
$ws.Range("F8").Value = 30.05354142189026
$ws.Range("F9").Value = 30.03009390830994
$ws.Range("F10").Value = 29.9560387134552
$ws.Range("F11").Value = 30.39167523384094
$ws.Range("F12").Value = 29.47680115699768
$ws.Range("F13").Value = 29.51582789421081
$ws.Range("F14").Value = 29.57950401306152
$ws.Range("F15").Value = 29.52680897712708
$ws.Range("F16").Value = 30.77045845985413
$ws.Range("F17").Value = 30.95502185821533
$ws.Range("F18").Value = 30.33430099487305
$ws.Range("F19").Value = 29.94668078422546
$ws.Range("F20").Value = 29.98090624809265
$ws.Range("F21").Value = 30.32273626327514

$ws = $wb.Worksheets.Item("run_4")
$ws.Range("F2").Value = 30.35609984397888
$ws.Range("F3").Value = 30.00747847557068
$ws.Range("F4").Value = 29.95002031326294
$ws.Range("F5").Value = 29.98706316947937
$ws.Range("F6").Value = 29.98426175117493
$ws.Range("F7").Value = 29.98045587539673
$ws.Range("F8").Value = 30.20563960075378
$ws.Range("F9").Value = 30.16924118995667
$ws.Range("F10").Value = 29.94273281097412
$ws.Range("F11").Value = 30.35927248001098
$ws.Range("F12").Value = 29.82886409759521
$ws.Range("F13").Value = 29.45768189430237
$ws.Range("F14").Value = 29.55277442932129
$ws.Range("F15").Value = 29.50424647331237
$ws.Range("F16").Value = 30.45827007293701
$ws.Range("F17").Value = 31.10278964042664
$ws.Range("F18").Value = 30.56849694252014
$ws.Range("F19").Value = 30.90485525131226
$ws.Range("F20").Value = 29.83219718933105
$ws.Range("F21").Value = 30.13261842727661

$ws = $wb.Worksheets.Item("run_5")
$ws.Range("F2").Value = 30.00175428390503
$ws.Range("F3").Value = 29.77335023880005
$ws.Range("F4").Value = 29.81047391891479
$ws.Range("F5").Value = 29.70942378044128
$ws.Range("F6").Value = 29.80865526199341
$ws.Range("F7").Value = 29.77254438400269
$ws.Range("F8").Value = 29.65243864059448
$ws.Range("F9").Value = 29.68417978286743
$ws.Range("F10").Value = 29.71243977546692
$ws.Range("F11").Value = 29.99475908279419
$ws.Range("F12").Value = 29.86923575401306
$ws.Range("F13").Value = 29.85220336914062
$ws.Range("F14").Value = 29.76237630844116
$ws.Range("F15").Value = 29.76931834220886
$ws.Range("F16").Value = 29.8070182800293
$ws.Range("F17").Value = 29.87999200820923
$ws.Range("F18").Value = 29.78793907165528
$ws.Range("F19").Value = 30.07143807411194
$ws.Range("F20").Value = 29.79715180397034
$ws.Range("F21").Value = 30.08678841590881
